$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 20:01"

$ws.Range("B4").Value = 5768900
$ws.Range("C4").Value = 22628
$ws.Range("D4").Value = 3103348
$ws.Range("E4").Value = 2487721
$ws.Range("G4").Value = 407
$ws.Range("H4").Value = 177831

$ws.Range("B6").Value = 2973317
$ws.Range("C6").Value = 68988
$ws.Range("D6").Value = 2218306
$ws.Range("E6").Value = 699083
$ws.Range("G6").Value = 953
$ws.Range("H6").Value = 55928

$ws.Range("B13").Value = 393769
$ws.Range("C13").Value = 1920
$ws.Range("D13").Value = 367897
$ws.Range("E13").Value = 15149
$ws.Range("G13").Value = 52
$ws.Range("H13").Value = 10723

$ws.Range("B21").Value = 255723
$ws.Range("C21").Value = 1203
$ws.Range("D21").Value = 235569
$ws.Range("E21").Value = 14074
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 6080

$ws.Range("A22").Value = "Francia"
$ws.Range("B22").Value = 234400
$ws.Range("C22").Value = 4586
$ws.Range("D22").Value = 84642
$ws.Range("E22").Value = 119255
$ws.Range("G22").Value = 23
$ws.Range("H22").Value = 30503

$ws.Range("A23").Value = "Alemania"
$ws.Range("B23").Value = 231389
$ws.Range("C23").Value = 105
$ws.Range("D23").Value = 205800
$ws.Range("E23").Value = 16264
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 9325

$ws.Range("B27").Value = 124099
$ws.Range("C27").Value = 226
$ws.Range("D27").Value = 110484
$ws.Range("E27").Value = 4555
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 9060

$ws.Range("A29").Value = "Ecuador"
$ws.Range("B29").Value = 106481
$ws.Range("C29").Value = 973
$ws.Range("D29").Value = 87730
$ws.Range("E29").Value = 12503
$ws.Range("G29").Value = 48
$ws.Range("H29").Value = 6248

$ws.Range("A30").Value = "Bolivia"
$ws.Range("B30").Value = 106065
$ws.Range("C30").Value = 1015
$ws.Range("D30").Value = 41111
$ws.Range("E30").Value = 60649
$ws.Range("G30").Value = 72
$ws.Range("H30").Value = 4305

$ws.Range("B35").Value = 89867
$ws.Range("C35").Value = 857
$ws.Range("D35").Value = 59132
$ws.Range("E35").Value = 29202
$ws.Range("G35").Value = 28
$ws.Range("H35").Value = 1533

$ws.Range("A53").Value = "Marruecos"
$ws.Range("B53").Value = 49247
$ws.Range("C53").Value = 1609
$ws.Range("D53").Value = 34199
$ws.Range("E53").Value = 14231
$ws.Range("G53").Value = 42
$ws.Range("H53").Value = 817

$ws.Range("A54").Value = "Barein"
$ws.Range("B54").Value = 48303
$ws.Range("D54").Value = 44628
$ws.Range("E54").Value = 3494
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 181

$ws.Range("B58").Value = 40667
$ws.Range("C58").Value = 409
$ws.Range("D58").Value = 28587
$ws.Range("E58").Value = 10662
$ws.Range("G58").Value = 7
$ws.Range("H58").Value = 1418

$ws.Range("A62").Value = "Etiopia"
$ws.Range("B62").Value = 37665
$ws.Range("C62").Value = 1829
$ws.Range("D62").Value = 13913
$ws.Range("E62").Value = 23115
$ws.Range("G62").Value = 17
$ws.Range("H62").Value = 637

$ws.Range("A63").Value = "Venezuela"
$ws.Range("B63").Value = 37567
$ws.Range("D63").Value = 26330
$ws.Range("E63").Value = 10926
$ws.Range("H63").Value = 311

$ws.Range("B70").Value = 27755
$ws.Range("C70").Value = 79
$ws.Range("E70").Value = 2615

$ws.Range("B88").Value = 10627
$ws.Range("C88").Value = 255
$ws.Range("E88").Value = 1224
$ws.Range("G88").Value = 3
$ws.Range("H88").Value = 277

$ws.Range("B104").Value = 6564
$ws.Range("C104").Value = 194
$ws.Range("D104").Value = 4012
$ws.Range("E104").Value = 2527

$ws.Range("A108").Value = "Namibia"
$ws.Range("B108").Value = 5227
$ws.Range("C108").Value = 315
$ws.Range("D108").Value = 2457
$ws.Range("E108").Value = 2728
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 42

$ws.Range("A109").Value = "Hungria"
$ws.Range("B109").Value = 5098
$ws.Range("C109").Value = 52
$ws.Range("D109").Value = 3681
$ws.Range("E109").Value = 806
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 611

$ws.Range("B115").Value = 4128
$ws.Range("C115").Value = 18
$ws.Range("D115").Value = 2682
$ws.Range("E115").Value = 1365

$ws.Range("B125").Value = 2941
$ws.Range("C125").Value = 23
$ws.Range("E125").Value = 141

$ws.Range("B139").Value = 1972
$ws.Range("C139").Value = 3
$ws.Range("D139").Value = 1542
$ws.Range("E139").Value = 361

$ws.Range("A157").Value = "Reunion"
$ws.Range("B157").Value = 1075
$ws.Range("C157").Value = 79
$ws.Range("D157").Value = 692
$ws.Range("E157").Value = 378
$ws.Range("H157").Value = 5

$ws.Range("A158").Value = "Principado de Andorra"
$ws.Range("B158").Value = 1045
$ws.Range("C158").Value = 21
$ws.Range("D158").Value = 875
$ws.Range("E158").Value = 117
$ws.Range("H158").Value = 53

$ws.Range("A159").Value = "Lesoto"
$ws.Range("B159").Value = 1015
$ws.Range("C159").Value = 19
$ws.Range("D159").Value = 472
$ws.Range("E159").Value = 513
$ws.Range("H159").Value = 30

$ws.Range("A160").Value = "Vietnam"
$ws.Range("B160").Value = 1009
$ws.Range("C160").Value = 2
$ws.Range("D160").Value = 545
$ws.Range("E160").Value = 439
$ws.Range("H160").Value = 25

$ws.Range("A189").Value = "Monaco"
$ws.Range("B189").Value = 154
$ws.Range("C189").Value = 4
$ws.Range("D189").Value = 116
$ws.Range("E189").Value = 34
$ws.Range("H189").Value = 4

$ws.Range("A190").Value = "Butan"
$ws.Range("B190").Value = 153
$ws.Range("C190").Value = 3
$ws.Range("D190").Value = 108
$ws.Range("E190").Value = 45
$ws.Range("H190").Value = 0
